$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.937.96"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.891.99"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8244"
$ws.Range("E5").Value = "  +5.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.61"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3247"
$ws.Range("E8").Value = "  +5.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.52"
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07039"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08030"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7474"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "1.891.55"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.209"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.21"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "29.930.47"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.894"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.92"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007760"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "2.141.79"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.924"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1565"
$ws.Range("E25").Value = "  +19.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.98"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.196"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.088"
$ws.Range("E29").Value = "  +3.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.520"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.274"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05642"
$ws.Range("E33").Value = "  +7.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.074"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.272"
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7297"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.723"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.779"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4425"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.93"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.955"
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8428"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.876"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.582"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.59"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.710"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "991.64"
$ws.Range("E49").Value = "  +6.77%  "
$ws.Range("D50").Value = "2.040.56"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.93"
$ws.Range("E51").Value = "  -0.39%  "
